$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "IU/mL"
$ws.Range("C3").Value = "mIU/mL"
$ws.Range("B4").Value = "t"
